{"js": "// Remove the duplicated \"Presentaci\u00f3n del proyecto\" screenshot paragraph.\n// The document already shows this image once under \"Introducci\u00f3n\" (the\n// featured/intro capture); the copy that was re-inserted right after the\n// \"Presentaci\u00f3n del proyecto\" heading is a duplicate and must be removed,\n// leaving the heading immediately followed by the explanatory paragraph\n// (\"Este proyecto demuestra c\u00f3mo integrar un \u2026\").\n\nconst body = context.document.body;\n\n// Locate the \"Presentaci\u00f3n del proyecto\" Heading2 paragraph. It is unique\n// in the document, so a text search pins it reliably without depending on\n// a fixed paragraph index.\nconst headingResults = body.search(\"Presentaci\u00f3n del proyecto\", { matchCase: true });\nheadingResults.load(\"items\");\nawait context.sync();\n\nif (headingResults.items.length > 0) {\n  const headingParagraph = headingResults.items[0].paragraphs.getFirst();\n  headingParagraph.load(\"styleBuiltIn\");\n  await context.sync();\n\n  // Only proceed if it is really the Heading2 (defensive \u2014 avoids acting on\n  // a stray body-text occurrence of the same words).\n  if (headingParagraph.styleBuiltIn === Word.BuiltInStyleName.heading2) {\n    const imageParagraph = headingParagraph.getNextOrNullObject();\n    imageParagraph.load(\"isNullObject\");\n    await context.sync();\n\n    if (!imageParagraph.isNullObject) {\n      const pictures = imageParagraph.inlinePictures;\n      pictures.load(\"items\");\n      await context.sync();\n\n      // Confirm the paragraph right after the heading is the duplicated\n      // centered picture paragraph before deleting it.\n      if (pictures.items.length > 0) {\n        imageParagraph.delete();\n        await context.sync();\n      }\n    }\n  }\n}\n", "ps1": "# Remove the duplicated \"Presentaci\u00f3n del proyecto\" screenshot paragraph.\n# The document already shows this image once under \"Introducci\u00f3n\" (the\n# featured/intro capture); the copy that was re-inserted right after the\n# \"Presentaci\u00f3n del proyecto\" heading is a duplicate and must be removed,\n# leaving the heading immediately followed by the explanatory paragraph\n# (\"Este proyecto demuestra c\u00f3mo integrar un \u2026\").\n\n$d = $word.ActiveDocument\n\n# Locate the \"Presentaci\u00f3n del proyecto\" Heading2 paragraph via Find \u2014 it is\n# unique in the document, so this pins it reliably without depending on a\n# fixed paragraph index.\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"Presentaci\u00f3n del proyecto\"\n$find.Forward = $true\n$found = $find.Execute()\n\nif ($found) {\n    # Resolve the paragraph that contains the Find hit.\n    $count = $d.Paragraphs.Count\n    $headingIndex = -1\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Start -le $range.Start -and $p.Range.End -ge $range.End) {\n            $headingIndex = $i\n            break\n        }\n    }\n\n    if ($headingIndex -gt 0 -and $headingIndex -lt $d.Paragraphs.Count) {\n        $imagePara = $d.Paragraphs.Item($headingIndex + 1)\n\n        # Only delete it if it is really the duplicated centered picture\n        # paragraph (defensive \u2014 avoids acting on unexpected content).\n        if ($imagePara.Range.InlineShapes.Count -gt 0) {\n            $imagePara.Range.Delete()\n        }\n    }\n}\n"}
